$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.645.47"
$ws.Range("E2").Value = "  -0.12%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.634.04"
$ws.Range("E3").Value = "  -0.32%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'212.21"
$ws.Range("C5").Copy()
$ws.Range("D5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E5").Value = "  -0.17%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.21%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.02%  "

# Row 8 - Solana
$ws.Range("D8").Value = "'23.30"
$ws.Range("C8").Copy()
$ws.Range("D8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E8").Value = "  +1.02%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +2.88%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.20%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.0856"
$ws.Range("C11").Copy()
$ws.Range("D11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E11").Value = "  -4.13%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.865.59"
$ws.Range("E12").Value = "  -0.33%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.628.09"
$ws.Range("E13").Value = "  -0.77%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.46%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "'0.553"
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E15").Value = "  -1.05%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "'65.17"
$ws.Range("C16").Copy()
$ws.Range("D16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E16").Value = "  +0.81%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.626.60"
$ws.Range("E17").Value = "  -0.15%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "'230.16"
$ws.Range("C18").Copy()
$ws.Range("D18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E18").Value = "  -0.15%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  -0.29%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "'7.57"
$ws.Range("C20").Copy()
$ws.Range("D20").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E20").Value = "  -2.06%  "

# Row 21 - Dai
$ws.Range("D21").Value = "'0.999"
$ws.Range("C21").Copy()
$ws.Range("D21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E21").Value = "  -0.02%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "'10.64"
$ws.Range("C22").Copy()
$ws.Range("D22").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E22").Value = "  +3.98%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  +1.12%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "'2.10"
$ws.Range("C24").Copy()
$ws.Range("D24").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E24").Value = "  +3.17%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'148.91"
$ws.Range("C25").Copy()
$ws.Range("D25").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E25").Value = "  -1.45%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  -1.13%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -0.27%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  +0.02%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  -0.69%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.00%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -0.68%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -1.07%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.478.67"
$ws.Range("E33").Value = "  +1.49%  "

# Row 34 - InternetComputer
$ws.Range("D34").Value = "'3.08"
$ws.Range("C34").Copy()
$ws.Range("D34").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E34").Value = "  -1.31%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -2.22%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -1.28%  "

# Row 37 - TrustWalletToken
$ws.Range("D37").Value = "'0.963"
$ws.Range("C37").Copy()
$ws.Range("D37").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E37").Value = "  +7.44%  "

# Row 38 - ARBITRUM
$ws.Range("D38").Value = "'0.878"
$ws.Range("C38").Copy()
$ws.Range("D38").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E38").Value = "  -0.07%  "

# Row 39 - ImmutableX
$ws.Range("D39").Value = "'0.558"
$ws.Range("C39").Copy()
$ws.Range("D39").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E39").Value = "  -1.30%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -0.02%  "

# Row 41 - WEMIXToken
$ws.Range("E41").Value = "  +1.42%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.02%  "

# Row 43 - Aave
$ws.Range("D43").Value = "'67.73"
$ws.Range("C43").Copy()
$ws.Range("D43").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E43").Value = "  -3.32%  "

# Row 44 - mCoin
$ws.Range("E44").Value = "  +0.32%  "

# Row 45 - MXToken
$ws.Range("E45").Value = "  -1.65%  "

# Row 46 - FraxShare
$ws.Range("E46").Value = "  -5.29%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.775.21"
$ws.Range("E47").Value = "  -0.38%  "

# Row 48 - RenderToken
$ws.Range("E48").Value = "  +0.05%  "

# Row 49 - Quant
$ws.Range("D49").Value = "'87.51"
$ws.Range("C49").Copy()
$ws.Range("D49").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E49").Value = "  +0.77%  "

# Row 50 - BabyDogeCoin
$ws.Range("E50").Value = "  -1.52%  "

# Row 51 - Algorand
$ws.Range("D51").Value = "'0.0991"
$ws.Range("C51").Copy()
$ws.Range("D51").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E51").Value = "  -0.09%  "

$excel.CutCopyMode = $false
